# Add new worksheet "all_possible_combinations" after the last existing sheet ("mapping")
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "all_possible_combinations"

# Header row
$ws.Cells.Item(1,1).Value = "v1"
$ws.Cells.Item(1,2).Value = "v2"
$ws.Cells.Item(1,3).Value = "v3"
$ws.Cells.Item(1,4).Value = "v4"

# All 24 permutations of (1,2,3,4) - the "all possible combinations" of a sudoku row
$ws.Cells.Item(2,1).Value = 2
$ws.Cells.Item(2,2).Value = 1
$ws.Cells.Item(2,3).Value = 3
$ws.Cells.Item(2,4).Value = 4
$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = 1
$ws.Cells.Item(3,3).Value = 4
$ws.Cells.Item(3,4).Value = 3
$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,2).Value = 1
$ws.Cells.Item(4,3).Value = 4
$ws.Cells.Item(4,4).Value = 2
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = 1
$ws.Cells.Item(5,3).Value = 2
$ws.Cells.Item(5,4).Value = 4
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = 1
$ws.Cells.Item(6,3).Value = 3
$ws.Cells.Item(6,4).Value = 2
$ws.Cells.Item(7,1).Value = 4
$ws.Cells.Item(7,2).Value = 1
$ws.Cells.Item(7,3).Value = 2
$ws.Cells.Item(7,4).Value = 3
$ws.Cells.Item(8,1).Value = 1
$ws.Cells.Item(8,2).Value = 2
$ws.Cells.Item(8,3).Value = 3
$ws.Cells.Item(8,4).Value = 4
$ws.Cells.Item(9,1).Value = 1
$ws.Cells.Item(9,2).Value = 2
$ws.Cells.Item(9,3).Value = 4
$ws.Cells.Item(9,4).Value = 3
$ws.Cells.Item(10,1).Value = 3
$ws.Cells.Item(10,2).Value = 2
$ws.Cells.Item(10,3).Value = 1
$ws.Cells.Item(10,4).Value = 4
$ws.Cells.Item(11,1).Value = 3
$ws.Cells.Item(11,2).Value = 2
$ws.Cells.Item(11,3).Value = 4
$ws.Cells.Item(11,4).Value = 1
$ws.Cells.Item(12,1).Value = 4
$ws.Cells.Item(12,2).Value = 2
$ws.Cells.Item(12,3).Value = 1
$ws.Cells.Item(12,4).Value = 3
$ws.Cells.Item(13,1).Value = 4
$ws.Cells.Item(13,2).Value = 2
$ws.Cells.Item(13,3).Value = 3
$ws.Cells.Item(13,4).Value = 1
$ws.Cells.Item(14,1).Value = 1
$ws.Cells.Item(14,2).Value = 3
$ws.Cells.Item(14,3).Value = 4
$ws.Cells.Item(14,4).Value = 2
$ws.Cells.Item(15,1).Value = 1
$ws.Cells.Item(15,2).Value = 3
$ws.Cells.Item(15,3).Value = 2
$ws.Cells.Item(15,4).Value = 4
$ws.Cells.Item(16,1).Value = 2
$ws.Cells.Item(16,2).Value = 3
$ws.Cells.Item(16,3).Value = 4
$ws.Cells.Item(16,4).Value = 1
$ws.Cells.Item(17,1).Value = 2
$ws.Cells.Item(17,2).Value = 3
$ws.Cells.Item(17,3).Value = 1
$ws.Cells.Item(17,4).Value = 4
$ws.Cells.Item(18,1).Value = 4
$ws.Cells.Item(18,2).Value = 3
$ws.Cells.Item(18,3).Value = 2
$ws.Cells.Item(18,4).Value = 1
$ws.Cells.Item(19,1).Value = 4
$ws.Cells.Item(19,2).Value = 3
$ws.Cells.Item(19,3).Value = 1
$ws.Cells.Item(19,4).Value = 2
$ws.Cells.Item(20,1).Value = 1
$ws.Cells.Item(20,2).Value = 4
$ws.Cells.Item(20,3).Value = 3
$ws.Cells.Item(20,4).Value = 2
$ws.Cells.Item(21,1).Value = 1
$ws.Cells.Item(21,2).Value = 4
$ws.Cells.Item(21,3).Value = 2
$ws.Cells.Item(21,4).Value = 3
$ws.Cells.Item(22,1).Value = 2
$ws.Cells.Item(22,2).Value = 4
$ws.Cells.Item(22,3).Value = 3
$ws.Cells.Item(22,4).Value = 1
$ws.Cells.Item(23,1).Value = 2
$ws.Cells.Item(23,2).Value = 4
$ws.Cells.Item(23,3).Value = 1
$ws.Cells.Item(23,4).Value = 3
$ws.Cells.Item(24,1).Value = 3
$ws.Cells.Item(24,2).Value = 4
$ws.Cells.Item(24,3).Value = 2
$ws.Cells.Item(24,4).Value = 1
$ws.Cells.Item(25,1).Value = 3
$ws.Cells.Item(25,2).Value = 4
$ws.Cells.Item(25,3).Value = 1
$ws.Cells.Item(25,4).Value = 2

# Turn the range into an Excel Table (ListObject) named Table1, sorted by v2
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:D25"), [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table1"
$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($ws.Range("B1:B25"))
$lo.Sort.Header = 1
$lo.Sort.Apply()

# Selection/view state for the new sheet
[void]$ws.Range("C2").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
